# Apply the "#5: insurance, claim, debt, investment done" edit.
#
# Sheets "保險" (insurance, sheet6) and "債務" (debt, sheet7) are missing
# the trailing metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index - and for the debt
# sheet also register_date/register_reason/total/owner/debtor/species in
# the right slots) that every other sheet in this workbook already
# carries. Row 1 (the header row) also currently just echoes column B's
# data instead of real field names. This script fixes both sheets to
# match the standard layout already used by the other sheets
# (現金/存款/...).

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

function Set-WithFormat {
    # Writes $value into ($row, $col) of $ws, then stamps it with the
    # formatting of ($formatSourceRow, $formatSourceCol) so the new cell
    # matches the look of its neighbours (bold/bordered header cells vs.
    # plain body cells).
    param($ws, $row, $col, $value, $formatSourceRow, $formatSourceCol)

    $cell = $ws.Cells.Item($row, $col)

    # A handful of values look like ISO dates ("2013-07-11") and Excel
    # will silently reinterpret them as date serials on assignment. Force
    # the cell to Text first so the literal string is preserved, matching
    # the source data (stored as a plain shared string, not a date).
    if ($value -is [string] -and $value -match '^\d{4}-\d{2}-\d{2}$') {
        $cell.NumberFormat = "@"
    }

    $cell.Value = $value

    $ws.Cells.Item($formatSourceRow, $formatSourceCol).Copy()
    $cell.PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# Sheet 6: 保險 (insurance)
# ---------------------------------------------------------------------
$wsIns = $wb.Worksheets.Item(6)

# --- Row 1 (header) ---
Set-WithFormat $wsIns 1 2  "company"            1 2   # B1
Set-WithFormat $wsIns 1 3  "name"               1 2   # C1
Set-WithFormat $wsIns 1 4  "owner"              1 2   # D1
Set-WithFormat $wsIns 1 5  "property_category"  1 2   # E1
Set-WithFormat $wsIns 1 6  "category"           1 2   # F1
Set-WithFormat $wsIns 1 7  "date"               1 2   # G1
Set-WithFormat $wsIns 1 8  "legislator_name"    1 2   # H1
Set-WithFormat $wsIns 1 9  "legislator_id"      1 2   # I1
Set-WithFormat $wsIns 1 10 "source_file"        1 2   # J1
Set-WithFormat $wsIns 1 11 "index"              1 2   # K1

# --- Row 2 (data, index 82) ---
Set-WithFormat $wsIns 2 2  "富邦人壽"            2 2   # B2 (value unchanged)
Set-WithFormat $wsIns 2 3  "生存還本保險"         2 2   # C2 (value unchanged)
Set-WithFormat $wsIns 2 4  "林正二"              2 2   # D2 (value unchanged)
Set-WithFormat $wsIns 2 5  "insurance"          2 2   # E2
Set-WithFormat $wsIns 2 6  "normal"             2 2   # F2
Set-WithFormat $wsIns 2 7  "2013-07-11"         2 2   # G2
Set-WithFormat $wsIns 2 8  "林正二"              2 2   # H2
Set-WithFormat $wsIns 2 9  788                  2 2   # I2
Set-WithFormat $wsIns 2 10 "tmp685a1"           2 2   # J2
Set-WithFormat $wsIns 2 11 82                   2 2   # K2

# --- Row 3 (data, index 83) ---
Set-WithFormat $wsIns 3 2  "富邦人壽"            3 2   # B3 (value unchanged)
Set-WithFormat $wsIns 3 3  "年金保險"            3 2   # C3 (value unchanged)
Set-WithFormat $wsIns 3 4  "林正二"              3 2   # D3 (value unchanged)
Set-WithFormat $wsIns 3 5  "insurance"          3 2   # E3
Set-WithFormat $wsIns 3 6  "normal"             3 2   # F3
Set-WithFormat $wsIns 3 7  "2013-07-11"         3 2   # G3
Set-WithFormat $wsIns 3 8  "林正二"              3 2   # H3
Set-WithFormat $wsIns 3 9  788                  3 2   # I3
Set-WithFormat $wsIns 3 10 "tmp685a1"           3 2   # J3
Set-WithFormat $wsIns 3 11 83                   3 2   # K3

# ---------------------------------------------------------------------
# Sheet 7: 債務 (debt)
# ---------------------------------------------------------------------
$wsDebt = $wb.Worksheets.Item(7)

# --- Row 1 (header) ---
Set-WithFormat $wsDebt 1 2  "species"           1 2   # B1
Set-WithFormat $wsDebt 1 3  "debtor"            1 2   # C1
Set-WithFormat $wsDebt 1 4  "owner"             1 2   # D1
Set-WithFormat $wsDebt 1 5  "total"             1 2   # E1
Set-WithFormat $wsDebt 1 6  "register_date"     1 2   # F1
Set-WithFormat $wsDebt 1 7  "register_reason"   1 2   # G1
Set-WithFormat $wsDebt 1 8  "property_category" 1 2   # H1
Set-WithFormat $wsDebt 1 9  "category"          1 2   # I1
Set-WithFormat $wsDebt 1 10 "date"              1 2   # J1
Set-WithFormat $wsDebt 1 11 "legislator_name"   1 2   # K1
Set-WithFormat $wsDebt 1 12 "legislator_id"     1 2   # L1
Set-WithFormat $wsDebt 1 13 "source_file"       1 2   # M1
Set-WithFormat $wsDebt 1 14 "index"             1 2   # N1

# --- Row 2 (data, index 93) ---
Set-WithFormat $wsDebt 2 2  "房屋貸款"                                    2 2   # B2 (value unchanged)
Set-WithFormat $wsDebt 2 3  "林正二"                                     2 2   # C2
Set-WithFormat $wsDebt 2 4  "臺灣土地銀行台東分行臺東縣台東市中華路1段357號"    2 2   # D2
Set-WithFormat $wsDebt 2 5  159006                                     2 2   # E2
Set-WithFormat $wsDebt 2 6  "83年04月30日"                               2 2   # F2
Set-WithFormat $wsDebt 2 7  "房屋貸款"                                    2 2   # G2
Set-WithFormat $wsDebt 2 8  "debt"                                     2 2   # H2
Set-WithFormat $wsDebt 2 9  "normal"                                   2 2   # I2
Set-WithFormat $wsDebt 2 10 "2013-07-11"                                2 2   # J2
Set-WithFormat $wsDebt 2 11 "林正二"                                     2 2   # K2
Set-WithFormat $wsDebt 2 12 788                                         2 2   # L2
Set-WithFormat $wsDebt 2 13 "tmp685a1"                                  2 2   # M2
Set-WithFormat $wsDebt 2 14 93                                          2 2   # N2

# --- Row 3 (data, index 94) ---
Set-WithFormat $wsDebt 3 2  "房屋貸款"                                    3 2   # B3 (value unchanged)
Set-WithFormat $wsDebt 3 3  "林正二"                                     3 2   # C3
Set-WithFormat $wsDebt 3 4  "臺灣土地銀行花蓮分行花蓮縣花蓮市中山路356號"      3 2   # D3
Set-WithFormat $wsDebt 3 5  3408061                                    3 2   # E3
Set-WithFormat $wsDebt 3 6  "88年04月07日"                               3 2   # F3
Set-WithFormat $wsDebt 3 7  "房屋貸款"                                    3 2   # G3
Set-WithFormat $wsDebt 3 8  "debt"                                     3 2   # H3
Set-WithFormat $wsDebt 3 9  "normal"                                   3 2   # I3
Set-WithFormat $wsDebt 3 10 "2013-07-11"                                3 2   # J3
Set-WithFormat $wsDebt 3 11 "林正二"                                     3 2   # K3
Set-WithFormat $wsDebt 3 12 788                                         3 2   # L3
Set-WithFormat $wsDebt 3 13 "tmp685a1"                                  3 2   # M3
Set-WithFormat $wsDebt 3 14 94                                          3 2   # N3

# --- Row 4 (data, index 95) ---
Set-WithFormat $wsDebt 4 2  "房屋貸款"                                    4 2   # B4 (value unchanged)
Set-WithFormat $wsDebt 4 3  "林正二"                                     4 2   # C4
Set-WithFormat $wsDebt 4 4  "華南商業銀行花蓮分行花蓮縣花蓮市中山路78號"       4 2   # D4
Set-WithFormat $wsDebt 4 5  2541362                                    4 2   # E4
Set-WithFormat $wsDebt 4 6  "88年04月01日"                               4 2   # F4
Set-WithFormat $wsDebt 4 7  "房屋貸款"                                    4 2   # G4
Set-WithFormat $wsDebt 4 8  "debt"                                     4 2   # H4
Set-WithFormat $wsDebt 4 9  "normal"                                   4 2   # I4
Set-WithFormat $wsDebt 4 10 "2013-07-11"                                4 2   # J4
Set-WithFormat $wsDebt 4 11 "林正二"                                     4 2   # K4
Set-WithFormat $wsDebt 4 12 788                                         4 2   # L4
Set-WithFormat $wsDebt 4 13 "tmp685a1"                                  4 2   # M4
Set-WithFormat $wsDebt 4 14 95                                          4 2   # N4
